# Add three new ObjectLocator rows to the "Lease" sheet of the Common.xlsx
# object-locator workbook: codeType / contractType dropdown locators and a
# "selectContract" popup locator (supports: "delete multiple leases from a
# property", "set Log message from sheet", "verify panel content false").

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Lease")

# Row 19 - codeType dropdown listbox locators
$ws.Range("A19").Value2 = "codeType"
$ws.Range("B19").Value2 = "by_xpath"
$ws.Range("C19").Value2 = "//*[contains(@aria-owns,'CodeTypeID_listbox')]"
$ws.Range("D19").Value2 = "//*[contains(@id,'CodeTypeID_listbox')]"

# Row 20 - contractType dropdown listbox locators
$ws.Range("A20").Value2 = "contractType"
$ws.Range("B20").Value2 = "by_xpath"
$ws.Range("C20").Value2 = "//*[contains(@aria-owns,'ContractTypeID_listbox')]"
$ws.Range("D20").Value2 = "//*[contains(@id,'ContractTypeID_listbox')]"

# Row 21 - selectContract popup OK button locator
$ws.Range("A21").Value2 = "selectContract"
$ws.Range("B21").Value2 = "by_xpath"
$ws.Range("C21").Value2 = "//*[@id='selectContract_Ok']"

# Match the sheet's recorded selection/active cell after the edit
$ws.Range("A21").Select()
